# Apply CrossRef enrichment + author-list re-serialization updates
# as described by the commit "added and updated with crossref".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


# Row 2: new CrossRef match for 10.4269/ajtmh.19-0311
$ws.Range("C2").Value = '"Women\u2019s Health Perceptions and Beliefs Related to Zika Virus Exposure during the 2016 Outbreak in Northern Brazil"'
$ws.Range("E2").Value = '[Elizabeth J.%Anderson%xref no email%1,  Kacey C.%Ernst%xref no email%1,  Francisco Fernando%Martins%xref no email%1,  Cicera da Silva%Martins%xref no email%1,  Mary P.%Koss%xref no email%1]'
$ws.Range("F2").Value = '10.4269/ajtmh.19-0311'
$ws.Range("G2").Value = 'CROSSREF'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '2023-05-16'

# Row 3: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E3").Value = '[Ana Rosa%Linde-Arias%linde14@yahoo.com%1,          Maria%Roura%NULL%1,          Eduardo%Siqueira%NULL%1]'

# Row 4: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E4").Value = '[Maria%Belizan%mbelizan@iecs.org.ar%1,          Edna%Maradiaga%edjamar3006@yahoo.com%1,          Javier%Roberti%jroberti@iecs.org.ar%1,          Maricela%Casco-Aguilar%marykasco@yahoo.com%1,          Alison F.%Ortez%alison_fabiola@yahoo.es%1,          Juan C.%Avila-Flores%javilaflores3@gmail.com%1,          Gloria%González%marilyntoin@yahoo.com%1,          Carolina%Bustillo%mcbu1502@yahoo.com%0,          Alejandra%Calderón%lilianalecalderon@gmail.com%1,          Harry%Bock%hbockme@hotmail.com%1,          María L.%Cafferata%NULL%1,          Adriano B.%Tavares%adriano.b.tavares@gmail.com%1,          Jackeline%Alger%jackelinealger@gmail.com%1,          Moazzam%Ali%alimoa@who.int%1]'

# Row 5: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E5").Value = '[ Catherine A.%Boyd%null%1,           Julie A.%Gazmararian%null%1,           Winifred Wilkins%Thompson%null%1,         Catherine A.%Boyd%null%1,         Julie A.%Gazmararian%null%1,         Winifred Wilkins%Thompson%null%1]'

# Row 6: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E6").Value = '[Benner%P.%coreGivesNoEmail%1,         Bohme%G.%coreGivesNoEmail%1,         Chiang%H-H%coreGivesNoEmail%4,         Chiang%H-H%coreGivesNoEmail%0,         Chiang%H-H%coreGivesNoEmail%0,         Chiang%H-H%coreGivesNoEmail%0,         Donohoe%J.%coreGivesNoEmail%1,         Hsien-Hsien%Chiang%coreGivesNoEmail%1,         I-Ling%Sue%coreGivesNoEmail%1,         International%Council of Nurses%coreGivesNoEmail%1,         Mei-Bih%Chen%coreGivesNoEmail%1,         Morgan%DL%coreGivesNoEmail%1,         Packer%MJ%coreGivesNoEmail%1]'

# Row 7: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E7").Value = '[Chee%YO%coreGivesNoEmail%1,         Dodgson%JE%coreGivesNoEmail%1,         Tarrant%M%coreGivesNoEmail%1,         Watkins%A%coreGivesNoEmail%1]'

# Row 8: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E8").Value = '[ Elin%Erland%null%1,           Bente%Dahl%null%1,         Elin%Erland%null%1,         Bente%Dahl%null%1]'

# Row 9: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E9").Value = '[Davis%Mark%coreGivesNoEmail%1,         Flowers%Paul%coreGivesNoEmail%1,         Lohm%Davina%coreGivesNoEmail%1,         Stephenson%Niamh%coreGivesNoEmail%1,         Waller%Emily%coreGivesNoEmail%1]'

# Row 10: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E10").Value = '[Hector M.%Gomez%NULL%1,          Carlos%Mejia Arbelaez%NULL%1,          Jovana A.%Ocampo Cañas%ja.ocampo@uniandes.edu.co%1]'

# Row 11: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E11").Value = '[ Susan%Jones%null%1,           Betty%Sam%null%1,           Florence%Bull%null%1,           Steven Bagie%Pieh%null%1,           Jaki%Lambert%null%1,           Florence%Mgawadere%null%1,           Somasundari%Gopalakrishnan%null%1,           Charles A.%Ameh%null%1,           Nynke%van den Broek%null%1,         Susan%Jones%null%1,         Betty%Sam%null%1,         Florence%Bull%null%1,         Steven Bagie%Pieh%null%1,         Jaki%Lambert%null%1,         Florence%Mgawadere%null%1,         Somasundari%Gopalakrishnan%null%1,         Charles A.%Ameh%null%1,         Nynke%van den Broek%null%1]'

# Row 12: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E12").Value = '[Caroline S.E.%Homer%NULL%1,          Miranda%Davies-Tuck%NULL%1,          Hannah G.%Dahlen%NULL%1,          Vanessa L.%Scarf%NULL%1]'

# Row 13: new CrossRef match for 10.1590/0102-311X00178917
$ws.Range("C13").Value = '"Women\u2019s lives in times of Zika: mosquito-controlled lives?"'
$ws.Range("D13").Value = '"Zika virus infection during pregnancy is a cause of congenital brain abnormalities.
 Its consequences to pregnancies has made governments, national and international agencies issue advices and recommendations to women.
 There is a clear need to investigate how the Zika outbreak affects the decisions that women take concerning their lives and the life of their families, as well as how women are psychologically and emotionally dealing with the outbreak.
 We conducted a qualitative study to address the impact of the Zika epidemic on the family life of women living in Brazil, Puerto Rico, and the US, who were affected by it to shed light on the social repercussions of Zika.
 Women were recruited through the snowball sampling technique and data was collected through semi-structured interviews.
 We describe the effects in mental health and the coping strategies that women use to deal with the Zika epidemic.
 Zika is taking a heavy toll on women\u2019s emotional well-being.
 They are coping with feelings of fear, helplessness, and uncertainty by taking drastic precautions to avoid infection that affect all areas of their lives.
 Coping strategies pose obstacles in professional life, lead to social isolation, including from family and partner, and threaten the emotional and physical well-being of women.
 Our findings suggest that the impacts of the Zika epidemic on women may be universal and global.
 Zika infection is a silent and heavy burden on women\u2019s shoulders.
'
$ws.Range("E13").Value = '[Ana Rosa%Linde%xref no email%1,  Carlos Eduardo%Siqueira%xref no email%1]'
$ws.Range("F13").Value = '10.1590/0102-311X00178917'
$ws.Range("G13").Value = 'CROSSREF'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '2023-05-16'

# Row 14: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E14").Value = '[ Molly M.%Lynch%null%1,           Elizabeth W.%Mitchell%null%1,           Jennifer L.%Williams%null%1,           Kelly%Brumbaugh%null%1,           Michelle%Jones-Bell%null%1,           Debra E.%Pinkney%null%1,           Christine M.%Layton%null%1,           Patricia W.%Mersereau%null%1,           Juliette S.%Kendrick%null%1,           Paula Eguino%Medina%null%1,           Lucia Rojas%Smith%null%1,         Molly M.%Lynch%null%1,         Elizabeth W.%Mitchell%null%1,         Jennifer L.%Williams%null%1,         Kelly%Brumbaugh%null%1,         Michelle%Jones-Bell%null%1,         Debra E.%Pinkney%null%1,         Christine M.%Layton%null%1,         Patricia W.%Mersereau%null%1,         Juliette S.%Kendrick%null%1,         Paula Eguino%Medina%null%1,         Lucia Rojas%Smith%null%1]'

# Row 15: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E15").Value = '[Niuniu%Sun%NULL%1,          Luoqun%Wei%NULL%1,          Suling%Shi%NULL%1,          Dandan%Jiao%NULL%1,          Runluo%Song%NULL%1,          Lili%Ma%NULL%1,          Hongwei%Wang%NULL%1,          Chao%Wang%NULL%1,          Zhaoguo%Wang%NULL%1,          Yanli%You%NULL%1,          Shuhua%Liu%NULL%1,          Hongyun%Wang%NULL%1]'

# Row 16: new CrossRef match for 10.1080/07399332.2020.1833012
$ws.Range("C16").Value = '"Effect of COVID-19 pandemic on women\u2019s health and safety: A study of immigrant survivors of intimate partner violence"'
$ws.Range("E16").Value = '[Bushra%Sabri%xref no email%1,  Maria%Hartley%xref no email%1,  Jyoti%Saha%xref no email%1,  Sarah%Murray%xref no email%1,  Nancy%Glass%xref no email%1,  Jacquelyn C.%Campbell%xref no email%1]'
$ws.Range("F16").Value = '10.1080/07399332.2020.1833012'
$ws.Range("G16").Value = 'CROSSREF'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '2023-05-25'

# Row 17: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E17").Value = '[Veronika%Tirado%Veronika.Tirado@ki.se%1,          Santiago A.%Morales Mesa%NULL%1,          John%Kinsman%NULL%1,          Anna Mia%Ekström%NULL%1,          Berta N.%Restrepo Jaramillo%NULL%1]'

# Row 18: Authors re-serialized (CORE lookup formatting refresh)
$ws.Range("E18").Value = '[Simon N%Williams%NULL%1,          Christopher J%Armitage%NULL%2,          Christopher J%Armitage%NULL%0,          Tova%Tampe%NULL%1,          Kimberly%Dienes%NULL%1]'
